$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells (rows 2-29) ---
$ws.Cells.Item(2, 2).Value = "NSE:3IINFOLTD"
$ws.Cells.Item(2, 3).Value = "NSE:ADORWELD"
$ws.Cells.Item(2, 4).Value = "NSE:CGPOWER"
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(3, 2).Value = "NSE:ACCURACY"
$ws.Cells.Item(3, 3).Value = "NSE:DHRUV"
$ws.Cells.Item(3, 4).Value = "NSE:ICICIGI"
$ws.Cells.Item(3, 6).Value = "NSE:GMRINFRA"
$ws.Cells.Item(4, 2).Value = "NSE:AHLEAST"
$ws.Cells.Item(4, 3).Value = "NSE:EPL"
$ws.Cells.Item(4, 4).Value = "NSE:SAIL"
$ws.Cells.Item(4, 6).Value = "NSE:NMDC"
$ws.Cells.Item(5, 2).Value = "NSE:AMBIKCO"
$ws.Cells.Item(5, 3).Value = "NSE:IDBI"
$ws.Cells.Item(6, 2).Value = "NSE:ARENTERP"
$ws.Cells.Item(6, 3).Value = "NSE:IOB"
$ws.Cells.Item(7, 2).Value = "NSE:BGRENERGY"
$ws.Cells.Item(7, 3).Value = "NSE:JAGSNPHARM"
$ws.Cells.Item(8, 2).Value = "NSE:BLS"
$ws.Cells.Item(8, 3).Value = "NSE:KIRLOSBROS"
$ws.Cells.Item(9, 2).Value = "NSE:BOMDYEING"
$ws.Cells.Item(9, 3).Value = "NSE:KOLTEPATIL"
$ws.Cells.Item(10, 2).Value = "NSE:CHALET"
$ws.Cells.Item(10, 3).Value = "NSE:NDTV"
$ws.Cells.Item(11, 2).Value = "NSE:CLSEL"
$ws.Cells.Item(11, 3).Value = "NSE:PEARLPOLY"
$ws.Cells.Item(12, 2).Value = "NSE:DLF"
$ws.Cells.Item(12, 3).Value = "NSE:RESPONIND"
$ws.Cells.Item(13, 2).Value = "NSE:DRCSYSTEMS"
$ws.Cells.Item(13, 3).Value = "NSE:RIIL"
$ws.Cells.Item(14, 2).Value = "NSE:GMRINFRA"
$ws.Cells.Item(14, 3).Value = ""
$ws.Cells.Item(15, 2).Value = "NSE:INDIANCARD"
$ws.Cells.Item(15, 3).Value = ""
$ws.Cells.Item(16, 2).Value = "NSE:INDNIPPON"
$ws.Cells.Item(17, 2).Value = "NSE:INDSWFTLTD"
$ws.Cells.Item(19, 2).Value = "NSE:KELLTONTEC"
$ws.Cells.Item(20, 2).Value = "NSE:LEMONTREE"
$ws.Cells.Item(21, 2).Value = "NSE:LICI"
$ws.Cells.Item(22, 2).Value = "NSE:LORDSCHLO"
$ws.Cells.Item(23, 2).Value = "NSE:LTFOODS"
$ws.Cells.Item(24, 2).Value = "NSE:MAHAPEXLTD"
$ws.Cells.Item(25, 2).Value = "NSE:MAXIND"
$ws.Cells.Item(26, 2).Value = "NSE:MIRCELECTR"
$ws.Cells.Item(27, 2).Value = "NSE:MTNL"
$ws.Cells.Item(28, 2).Value = "NSE:NAHARINDUS"
$ws.Cells.Item(29, 2).Value = "NSE:NEWGEN"

# --- Add new rows 30-33, copying style from row 2 (A column) for index cells ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(30, 1))
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(31, 1))
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(32, 1))
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(33, 1))
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "NSE:NLCINDIA"
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "NSE:NMDC"
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = "NSE:PKTEA"
$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 2).Value = "NSE:RVHL"
